$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-18 from 2023-09-09 (45178)
# to 2023-09-10 (45179), matching the serial-number bump in the source diff.
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
